$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D8").Value = "2016-03-10 04:58:42"
$wsZhCn.Range("G8").Value = "2016-03-10 04:59:26"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D8").Value = "2016-03-10 04:58:51"
$wsDeDe.Range("G8").Value = "2016-03-10 04:59:42"
